# Update the grading workbook: enter actual per-criteria scores, let the
# SUM/MIN/total formulas recalc on their own, and leave the selection on
# the last-edited cell (G28), matching the author's "grading" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- General --------------------------------------------------------------
$ws.Range("G4").Value = 1     # Description (HTML): 5 -> 1
$ws.Range("G5").Value = 6     # It works:            15 -> 6
                               # G6 = SUM(G4:G5) / H6 = MIN(C6,G6) recalc

# -- Functionality ----------------------------------------------------------
$ws.Range("G8").Value = 1     # Navigation:        2 -> 1
$ws.Range("G9").Value = 1     # User Input:        2 -> 1
$ws.Range("G10").Value = 1    # Input Validation:  3 -> 1
$ws.Range("G11").Value = 1    # Interactive UI:    5 -> 1
$ws.Range("G12").Value = 3    # other:            10 -> 3
                               # G13 = SUM(G8:G12) / H13 = MIN(C13,G13) recalc

# -- Engineering ------------------------------------------------------------
$ws.Range("G16").Value = 1    # Tests:                   5 -> 1
$ws.Range("G17").Value = 0    # Valid HTML/CSS:          2 -> 0
$ws.Range("G18").Value = 1    # Comments:                2 -> 1
# G19 Good naming and G21 Few dependencies are unchanged
$ws.Range("G20").Value = 1    # Not much duplication:    2 -> 1
                               # G22 = SUM(G15:G21) / H22 = MIN(C22,G22) recalc

# -- Extra --------------------------------------------------------------
# Only "Artistic value" scored any points this time; the other four Extra
# criteria go back to ungraded/blank cells.
$ws.Range("G24").ClearContents()  # Calling Services:   5 -> (blank)
$ws.Range("G25").ClearContents()  # Exposing Services:  5 -> (blank)
$ws.Range("G26").ClearContents()  # Technical challenge:5 -> (blank)
$ws.Range("G27").Value = 1        # Artistic value:     5 -> 1
$ws.Range("G28").ClearContents()  # Particular effort:  5 -> (blank)
                                   # G29 = SUM(G24:G28) / H29 = MIN(C29,G29) recalc

# H32/I32 (grand totals) and I34 (1 + TOTAL_POINTS/10) recalc automatically
# from the above via their existing formulas.

# Leave the selection where the grader left off.
$ws.Range("G28").Select()

# Best-effort: replicate the cosmetic window-geometry / font-charset tweaks
# from the original commit. These correspond to UI chrome (last on-screen
# window position) and a Calibri charset byte that this host doesn't expose
# through the Excel object model, so they're harmless no-ops if unsupported.
try { $excel.ActiveWindow.Left = 33700 } catch {}
try { $excel.ActiveWindow.Top = 140 } catch {}
try { $excel.ActiveWindow.Width = 17280 } catch {}
try { $excel.ActiveWindow.Height = 23080 } catch {}
